$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (pushes existing rows 2..26 down to 3..27).
# The inserted row borrows formatting from the row above (the bold header),
# so explicitly clear it to match the unstyled data rows used everywhere
# else in the sheet.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row with the new IPO record (에스엘에스바이오,
# listed 2023-10-20). Date-shaped values in columns A/O/P are stored as
# plain text (shared strings) throughout this sheet, not real dates, so a
# leading apostrophe forces text entry; ClearFormats afterwards drops the
# transient "quote prefix" style so the cell ends up with no style index,
# matching the rest of the column.
$ws.Range("A2").Value = "'2023-10-20"
$ws.Range("A2").ClearFormats()
$ws.Range("B2").Value = "에스엘에스바이오"
$ws.Range("C2").Value = "코스닥"
$ws.Range("D2").Value = 53.9
$ws.Range("E2").Value = "하나"
$ws.Range("F2").Value = 53.9
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "대표"
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = 7000
$ws.Range("N2").Value = 100
$ws.Range("O2").Value = "'2023-10-10"
$ws.Range("O2").ClearFormats()
$ws.Range("P2").Value = "'2023-10-13"
$ws.Range("P2").ClearFormats()
$ws.Range("Q2").Value = 577500
